$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest metric reading as a new row after the existing data.
$ws.Range("A73").Value = "2025-04-29 10:55:05"
$ws.Range("B73").Value = 214
